$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.496.26"
$ws.Range("E2").Value = "  -1.42%  "
$ws.Range("D3").Value = "3.412.39"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "3.417.20"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.124"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.442"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "3.995.33"
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000188"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.68%  "
$ws.Range("D17").Value = "64.520.79"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").Value = "3.387.65"
$ws.Range("E18").Value = "  -4.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.548"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000119"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.178"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0761"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.36%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "2.888.45"
$ws.Range("E40").Value = "  -4.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0317"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.771"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "321.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.47%  "
$ws.Range("E48").Value = "  +2.94%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.79%  "
